$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix bug in hookes_law: wrong velocity values (Damping coefficient row, column D)
$ws.Range("D3").Value = 0.05
$ws.Range("D4").Value = 0.05

# Update selection to D4
$ws.Range("D4").Select()
